$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.032.28"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").Value = "3.132.83"
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'587.35"
$ws.Range("E5").Value = "  -0.54%  "

$ws.Range("D6").Value = "'146.12"
$ws.Range("E6").Value = "  -0.44%  "

$ws.Range("D7").Value = "'1.01"
$ws.Range("E7").Value = "  +0.51%  "

$ws.Range("D8").Value = "3.129.69"
$ws.Range("E8").Value = "  +0.18%  "

$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -1.59%  "

$ws.Range("D10").Value = "'0.160"
$ws.Range("E10").Value = "  -2.54%  "

$ws.Range("D11").Value = "'5.85"
$ws.Range("E11").Value = "  +2.75%  "

$ws.Range("D12").Value = "'0.458"
$ws.Range("E12").Value = "  -2.06%  "

$ws.Range("D13").Value = "'0.0000247"
$ws.Range("E13").Value = "  -3.71%  "

$ws.Range("D14").Value = "'37.24"
$ws.Range("E14").Value = "  +3.43%  "

$ws.Range("D15").Value = "3.932.52"
$ws.Range("E15").Value = "  +7.77%  "

$ws.Range("E16").Value = "  -1.45%  "

$ws.Range("D17").Value = "3.184.91"
$ws.Range("E17").Value = "  +1.74%  "

$ws.Range("D18").Value = "63.918.80"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("D19").Value = "'7.06"
$ws.Range("E19").Value = "  -1.54%  "

$ws.Range("D20").Value = "'463.64"
$ws.Range("E20").Value = "  -0.55%  "

$ws.Range("D21").Value = "'14.32"
$ws.Range("E21").Value = "  +0.73%  "

$ws.Range("D22").Value = "'0.728"
$ws.Range("E22").Value = "  -0.57%  "

$ws.Range("D23").Value = "'7.39"
$ws.Range("E23").Value = "  -1.91%  "

$ws.Range("D24").Value = "'12.89"
$ws.Range("E24").Value = "  -2.89%  "

$ws.Range("D25").Value = "'80.84"
$ws.Range("E25").Value = "  -1.83%  "

$ws.Range("D26").Value = "'2.24"
$ws.Range("E26").Value = "  +3.99%  "

$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("D28").Value = "'9.32"
$ws.Range("E28").Value = "  +8.37%  "

$ws.Range("E29").Value = "  -1.17%  "

$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.19%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'2.20"
$ws.Range("E31").Value = "  -0.23%  "

$ws.Range("D32").Value = "'7.13"
$ws.Range("E32").Value = "  +4.06%  "

$ws.Range("D33").Value = "'26.91"
$ws.Range("E33").Value = "  -0.42%  "

$ws.Range("E34").Value = "  +0.29%  "

$ws.Range("D35").Value = "0.0₃0854"
$ws.Range("E35").Value = "  -1.86%  "

$ws.Range("D36").Value = "'1.04"
$ws.Range("E36").Value = "  -0.54%  "

$ws.Range("D37").Value = "'2.31"
$ws.Range("E37").Value = "  -3.46%  "

$ws.Range("E38").Value = "  -1.02%  "

$ws.Range("D39").Value = "'6.00"
$ws.Range("E39").Value = "  -2.09%  "

$ws.Range("D40").Value = "'51.43"
$ws.Range("E40").Value = "  +1.21%  "

$ws.Range("D41").Value = "'437.96"
$ws.Range("E41").Value = "  -2.42%  "

$ws.Range("E42").Value = "  +1.50%  "

$ws.Range("D43").Value = "'0.287"
$ws.Range("E43").Value = "  +3.34%  "

$ws.Range("D44").Value = "'0.0371"
$ws.Range("E44").Value = "  -0.48%  "

$ws.Range("D45").Value = "2.911.60"
$ws.Range("E45").Value = "  -0.29%  "

$ws.Range("D46").Value = "'39.70"
$ws.Range("E46").Value = "  +16.31%  "

$ws.Range("D47").Value = "'0.107"
$ws.Range("E47").Value = "  -3.58%  "

$ws.Range("D48").Value = "'126.68"
$ws.Range("E48").Value = "  -0.43%  "

$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("E50").Value = "  -1.09%  "

$ws.Range("D51").Value = "'2.19"
$ws.Range("E51").Value = "  +0.41%  "
